$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B17").Value = 5
$ws.Range("C17").Value = "Maximum Sum With Exactly K Elements"
$ws.Range("D17").Value = "Bosscoder Academy"

$ws.Range("C18").Value = "Richest Customer Wealth"
$ws.Range("D18").Value = "Bosscoder Academy"

$ws.Range("C19").Value = "1 to N (Recursion)"
$ws.Range("D19").Value = "Bosscoder Academy"

$ws.Range("C4").Select()
